$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refine the timestamp recorded for the existing last row (A54)
$ws.Range("A54").Value = 44367.76684843056

# Append the newly retrieved data row (row 55)
$ws.Range("A55").Value = 44368.76851705954
$ws.Range("B55").Value = 77979
$ws.Range("C55").Value = 65615
$ws.Range("D55").Value = 3467
$ws.Range("E55").Value = 2084
$ws.Range("F55").Value = 1480
$ws.Range("G55").Value = 20744
$ws.Range("H55").Value = 1411
$ws.Range("I55").Value = 884
$ws.Range("J55").Value = 178
